$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 162, pushing existing rows 162.. down by two.
$ws.Rows("162:163").Insert()

# New row 162: Feria Lagunitas de Puerto Montt, Apio, Primera, week of 2022-08-16 (serial 44789)
$ws.Range("A162").Value = 4
$ws.Range("B162").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C162").Value = "Los Lagos"
$ws.Range("D162").Value = 44789
$ws.Range("E162").Value = 10
$ws.Range("F162").Value = 100112017
$ws.Range("G162").Value = "Apio"
$ws.Range("H162").Value = "Americana (o)"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 25
$ws.Range("K162").Value = 14000
$ws.Range("L162").Value = 14000
$ws.Range("M162").Value = 14000
$ws.Range("N162").Value = "`$/docena de matas"
$ws.Range("O162").Value = "Región de Coquimbo"
$ws.Range("P162").Value = 2333
$ws.Range("Q162").Value = 6
$ws.Range("R162").Value = "Hortaliza"

# New row 163: same market/date, Segunda quality
$ws.Range("A163").Value = 4
$ws.Range("B163").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C163").Value = "Los Lagos"
$ws.Range("D163").Value = 44789
$ws.Range("E163").Value = 10
$ws.Range("F163").Value = 100112017
$ws.Range("G163").Value = "Apio"
$ws.Range("H163").Value = "Americana (o)"
$ws.Range("I163").Value = "Segunda"
$ws.Range("J163").Value = 25
$ws.Range("K163").Value = 12000
$ws.Range("L163").Value = 12000
$ws.Range("M163").Value = 12000
$ws.Range("N163").Value = "`$/docena de matas"
$ws.Range("O163").Value = "Región de Coquimbo"
$ws.Range("P163").Value = 2000
$ws.Range("Q163").Value = 6
$ws.Range("R163").Value = "Hortaliza"
